$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy'
$ws.Range("G3").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda'
$ws.Range("G4").Value = 'Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G5").Value = 'Dr. Veronia Rafat, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G6").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G7").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G8").Value = 'Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator'
$ws.Range("G9").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G10").Value = 'Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali'
$ws.Range("G11").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G12").Value = 'Dr. Salma El-Gendy, Administrator'
$ws.Range("G13").Value = 'D Wessam Atef, Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany'
$ws.Range("G17").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa'
$ws.Range("G23").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G24").Value = 'Dr. Salma Hassan, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Remon, Dr. Yasmin, Dr. Monica'
$ws.Range("G25").Value = 'Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Marina Atef'
$ws.Range("G26").Value = 'Dr. Youstina Magdy, Dr. Gehad Salah'
$ws.Range("G27").Value = 'Dr. Salma Hassan, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Remon, Dr. Yasmin'
$ws.Range("G28").Value = 'Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Remon, Dr. Nardine'
$ws.Range("G29").Value = 'Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica'
$ws.Range("G30").Value = 'Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G31").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda'
$ws.Range("G32").Value = 'Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G33").Value = 'Dr. Veronia Rafat, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G34").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G35").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G36").Value = 'Dr. Manar Montaser, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Asmaa Reda, Administrator'
$ws.Range("G37").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud'
$ws.Range("G38").Value = 'Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali'
$ws.Range("G39").Value = 'Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G40").Value = 'Dr. Salma El-Gendy, Administrator'
$ws.Range("G41").Value = 'D Wessam Atef, Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf, Dr. Omnia Mohammad, Dr. Safa Hany'
$ws.Range("G45").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Basma Hamed, Dr. Esraa Mostafa'
$ws.Range("G51").Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range("G52").Value = 'Dr. Salma Hassan, Dr. Aya Emad, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Remon, Dr. Yasmin, Dr. Monica'
$ws.Range("G53").Value = 'Dr. Aya Emad, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Remon, Dr. Marina Atef'
$ws.Range("G54").Value = 'Dr. Youstina Magdy, Dr. Gehad Salah'
$ws.Range("G55").Value = 'Dr. Salma Hassan, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Remon, Dr. Yasmin'
$ws.Range("G56").Value = 'Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Aya Hanafy, Dr. Neveen Nashaat, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Remon, Dr. Nardine'
$ws.Range("G57").Value = 'Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Remon, Dr. Monica'
